$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# abnormal.xlsx: "Enemy" column values (E2, E3) renamed to "target"
# (commit: enemy 삭제 / replaced by target column values)
$ws.Range("E2").Value2 = "target"
$ws.Range("E3").Value2 = "target"

# Update the UI selection state left behind by the edit (drag-selected
# J13:K13 with K13 as the active/last cell).
$ws.Range("J13:K13").Select() | Out-Null
